# Updated symbol list (prices / 1h volume % / a few reordered coin rows)
# as refreshed by the scheduled GitHub Actions scraper run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the numeric-looking Price/Volume cells as Text so Excel
# keeps the exact literal strings (matching the original inlineStr cells)
# instead of silently coercing them to numbers/percentages.
$textCells = @("D2", "E2", "D3", "E3", "D4", "E4", "E5", "D6", "E6", "D7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "E22", "D23", "E23", "D24", "E24", "E25", "D26", "E26", "E27", "E28", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "E46", "D47", "E47", "E48", "E49", "E50")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply the updated values cell by cell.
$ws.Range("D2").Value = "246.52"
$ws.Range("E2").Value = "0.53%"
$ws.Range("D3").Value = "29.74"
$ws.Range("E3").Value = "9.87%"
$ws.Range("D4").Value = "5.169"
$ws.Range("E4").Value = "1.70%"
$ws.Range("E5").Value = "0.34%"
$ws.Range("D6").Value = "6.600"
$ws.Range("E6").Value = "1.35%"
$ws.Range("D7").Value = "0.8559"
$ws.Range("D8").Value = "0.8692"
$ws.Range("E8").Value = "-1.18%"
$ws.Range("D9").Value = "0.1364"
$ws.Range("E9").Value = "2.82%"
$ws.Range("D10").Value = "0.07082"
$ws.Range("E10").Value = "2.49%"
$ws.Range("D11").Value = "0.02925"
$ws.Range("E11").Value = "3.40%"
$ws.Range("E12").Value = "-0.05%"
$ws.Range("D13").Value = "0.001511"
$ws.Range("E13").Value = "-0.40%"
$ws.Range("D14").Value = "0.04170"
$ws.Range("E14").Value = "2.02%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "0.006076"
$ws.Range("E15").Value = "-0.38%"
$ws.Range("B16").Value = "UpBots"
$ws.Range("C16").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D16").Value = "0.007489"
$ws.Range("E16").Value = "5,072.31%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.488"
$ws.Range("E17").Value = "-0.57%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "3.088"
$ws.Range("E18").Value = "2.60%"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "2.277"
$ws.Range("E19").Value = "2.16%"
$ws.Range("B20").Value = "One"
$ws.Range("C20").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D20").Value = "0.0006030"
$ws.Range("E20").Value = "0.73%"
$ws.Range("D21").Value = "0.3183"
$ws.Range("E21").Value = "0.20%"
$ws.Range("E22").Value = "5.53%"
$ws.Range("D23").Value = "0.1303"
$ws.Range("E23").Value = "2.28%"
$ws.Range("D24").Value = "3.468"
$ws.Range("E24").Value = "-2.40%"
$ws.Range("E25").Value = "0.49%"
$ws.Range("D26").Value = "0.005020"
$ws.Range("E26").Value = "26.31%"
$ws.Range("E27").Value = "0.29%"
$ws.Range("E28").Value = "22.27%"
$ws.Range("D40").Value = "0.03746"
$ws.Range("E40").Value = "0.71%"
$ws.Range("D41").Value = "0.005760"
$ws.Range("E41").Value = "67.35%"
$ws.Range("D42").Value = "0.1072"
$ws.Range("E42").Value = "1.40%"
$ws.Range("D43").Value = "0.002000"
$ws.Range("E43").Value = "-15.52%"
$ws.Range("D44").Value = "0.008327"
$ws.Range("E44").Value = "-11.27%"
$ws.Range("D45").Value = "0.00005212"
$ws.Range("E45").Value = "1.39%"
$ws.Range("E46").Value = "0.06%"
$ws.Range("D47").Value = "0.06470"
$ws.Range("E47").Value = "-36.22%"
$ws.Range("E48").Value = "-1.55%"
$ws.Range("E49").Value = "0.06%"
$ws.Range("E50").Value = "0.06%"
